$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.5271936044704619, 0.5352467862865621, 0.3759612643965448, 0.6131568024547593, 0.3249135209474168, 14),
    @(0.3608722674935116, 0.4023274030500674, 0.2193608715941454, 0.4683597672667298, 0.310740649621604, 13),
    @(0.2850699742386379, 0.3326698397266342, 0.1632228408311082, 0.4040084662864236, 0.2990128315009373, 12),
    @(0.437101989811909, 0.4532318499641988, 0.2559434172938412, 0.5059085068407539, 0.2671587441502564, 11),
    @(0.3896036158490621, 0.4137899396662769, 0.213010691573719, 0.4615308132440552, 0.2608100928752895, 10),
    @(0.3372085974095211, 0.3720251510852938, 0.1694528389792911, 0.4116464975914299, 0.2504218459199277, 9),
    @(0.3687366187234252, 0.3981931578693139, 0.1958236157540623, 0.4425196218859253, 0.2615490596241174, 8),
    @(0.3640834566367565, 0.3961645473499545, 0.1981723784143557, 0.4451655629250265, 0.2766795092775837, 7),
    @(0.3225455632229594, 0.3512104856557821, 0.1570729010663444, 0.3963242372935882, 0.2522790376821546, 6),
    @(0.3733140426150657, 0.3962272652827047, 0.198344628268327, 0.4453589880852603, 0.2715263657886741, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    for ($j = 0; $j -lt 6; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}
